# The original sheet had two header rows: row 1 held unit labels
# (feet, rho/rho_sl, P/P_sl, ...) and row 2 held the variable names
# (alt, sigma, delta, ...). This edit removes the units header row
# entirely so the variable-name row becomes the new row 1, and
# converts the altitude column from feet-by-10000 increments to
# feet-by-1000 increments (divide column A data by 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 1 (units row: feet, rho/rho_sl, P/P_sl, T/T_sl,
# Rankine, lb/ft^2, slugs/ft^3, ft/s, slugs/ft*s, ft^2/s, a/k.visc).
# This shifts every remaining row up by one, so the former row 2
# (alt, sigma, delta, theta, temp, press, dens, a, visc, k.visc,
# ratio) becomes the new row 1, and the data rows become rows 2-68.
$ws.Rows.Item(1).Delete()

# Rescale the altitude column (A) data rows from increments of
# 10000 ft down to increments of 1000 ft (divide by 10).
for ($r = 2; $r -le 68; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value() / 10
}

# Match the author's final selection.
[void]$ws.Range("K21").Select()
